$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header columns (row 1) ---
$ws.Cells.Item(1, 1).Value = 'mx_state'
$ws.Cells.Item(1, 2).Value = 'mx_municipality'
$ws.Cells.Item(1, 3).Value = 'n_matriculas'
$ws.Cells.Item(1, 4).Value = 'pct_matriculas'

# --- Title-case municipality/state names (Spanish connector words capitalised) ---
$ws.Cells.Item(7, 2).Value = 'Pabellón De Arteaga'
$ws.Cells.Item(8, 2).Value = 'Rincón De Romos'
$ws.Cells.Item(9, 2).Value = 'San Francisco De Los Romo'
$ws.Cells.Item(10, 2).Value = 'San José De Gracia'
$ws.Cells.Item(34, 2).Value = 'Amatenango De La Frontera'
$ws.Cells.Item(37, 2).Value = 'Bejucal De Ocampo'
$ws.Cells.Item(45, 2).Value = 'Chiapa De Corzo'
$ws.Cells.Item(50, 2).Value = 'Comitán De Domínguez'
$ws.Cells.Item(74, 2).Value = 'Marqués De Comillas'
$ws.Cells.Item(75, 2).Value = 'Mazapa De Madero'
$ws.Cells.Item(81, 2).Value = 'Ocozocoautla De Espinosa'
$ws.Cells.Item(89, 2).Value = 'Salto De Agua'
$ws.Cells.Item(90, 2).Value = 'San Cristóbal De Las Casas'
$ws.Cells.Item(129, 2).Value = 'Hidalgo Del Parral'
$ws.Cells.Item(141, 2).Value = 'San Francisco Del Oro'
$ws.Cells.Item(144, 2).Value = 'Valle De Zaragoza'
$ws.Cells.Item(169, 2).Value = 'San Juan De Sabinas'
$ws.Cells.Item(181, 2).Value = 'Villa De Álvarez'
$ws.Cells.Item(183, 1).Value = 'Ciudad De México'
$ws.Cells.Item(187, 2).Value = 'Cuajimalpa De Morelos'
$ws.Cells.Item(202, 2).Value = 'Coneto De Comonfort'
$ws.Cells.Item(216, 2).Value = 'Nombre De Dios'
$ws.Cells.Item(220, 2).Value = 'Pánuco De Coronado'
$ws.Cells.Item(226, 2).Value = 'San Juan De Guadalupe'
$ws.Cells.Item(227, 2).Value = 'San Juan Del Río'
$ws.Cells.Item(228, 2).Value = 'San Luis Del Cordero'
$ws.Cells.Item(229, 2).Value = 'San Pedro Del Gallo'
$ws.Cells.Item(239, 1).Value = 'Estado De México'
$ws.Cells.Item(239, 2).Value = 'Acambay De Ruíz Castañeda'
$ws.Cells.Item(242, 2).Value = 'Almoloya De Alquisiras'
$ws.Cells.Item(243, 2).Value = 'Almoloya De Juárez'
$ws.Cells.Item(250, 2).Value = 'Atizapán De Zaragoza'
$ws.Cells.Item(257, 2).Value = 'Coacalco De Berriozábal'
$ws.Cells.Item(263, 2).Value = 'Ecatepec De Morelos'
$ws.Cells.Item(270, 2).Value = 'Ixtapan De La Sal'
$ws.Cells.Item(271, 2).Value = 'Ixtapan Del Oro'
$ws.Cells.Item(284, 2).Value = 'Naucalpan De Juárez'
$ws.Cells.Item(293, 2).Value = 'San Felipe Del Progreso'
$ws.Cells.Item(295, 2).Value = 'San Simón De Guerrero'
$ws.Cells.Item(297, 2).Value = 'Soyaniquilpan De Juárez'
$ws.Cells.Item(306, 2).Value = 'Tenango Del Valle'
$ws.Cells.Item(317, 2).Value = 'Tlalnepantla De Baz'
$ws.Cells.Item(323, 2).Value = 'Valle De Bravo'
$ws.Cells.Item(324, 2).Value = 'Valle De Chalco Solidaridad'
$ws.Cells.Item(325, 2).Value = 'Villa De Allende'
$ws.Cells.Item(326, 2).Value = 'Villa Del Carbón'
$ws.Cells.Item(338, 2).Value = 'San Miguel De Allende'
$ws.Cells.Item(339, 2).Value = 'Apaseo El Alto'
$ws.Cells.Item(340, 2).Value = 'Apaseo El Grande'
$ws.Cells.Item(348, 2).Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Cells.Item(352, 2).Value = 'Jaral Del Progreso'
$ws.Cells.Item(360, 2).Value = 'Purísima Del Rincón'
$ws.Cells.Item(364, 2).Value = 'San Diego De La Unión'
$ws.Cells.Item(366, 2).Value = 'San Francisco Del Rincón'
$ws.Cells.Item(368, 2).Value = 'San Luis De La Paz'
$ws.Cells.Item(370, 2).Value = 'Santa Cruz De Juventino Rosas'
$ws.Cells.Item(372, 2).Value = 'Silao De La Victoria'
$ws.Cells.Item(377, 2).Value = 'Valle De Santiago'
$ws.Cells.Item(383, 2).Value = 'Acapulco De Juárez'
$ws.Cells.Item(385, 2).Value = 'Ajuchitlán Del Progreso'
$ws.Cells.Item(386, 2).Value = 'Alcozauca De Guerrero'
$ws.Cells.Item(390, 2).Value = 'Atenango Del Río'
$ws.Cells.Item(391, 2).Value = 'Atlamajalcingo Del Monte'
$ws.Cells.Item(393, 2).Value = 'Atoyac De Álvarez'
$ws.Cells.Item(394, 2).Value = 'Ayutla De Los Libres'
$ws.Cells.Item(397, 2).Value = 'Buenavista De Cuéllar'
$ws.Cells.Item(398, 2).Value = 'Chilapa De Álvarez'
$ws.Cells.Item(399, 2).Value = 'Chilpancingo De Los Bravo'
$ws.Cells.Item(400, 2).Value = 'Coahuayutla De José María Izazaga'
$ws.Cells.Item(405, 2).Value = 'Coyuca De Benítez'
$ws.Cells.Item(406, 2).Value = 'Coyuca De Catalán'
$ws.Cells.Item(410, 2).Value = 'Cuetzala Del Progreso'
$ws.Cells.Item(411, 2).Value = 'Cutzamala De Pinzón'
$ws.Cells.Item(417, 2).Value = 'Huitzuco De Los Figueroa'
$ws.Cells.Item(418, 2).Value = 'Iguala De La Independencia'
$ws.Cells.Item(420, 2).Value = 'Ixcateopan De Cuauhtémoc'
$ws.Cells.Item(423, 2).Value = 'La Unión De Isidoro Montes De Oca'
$ws.Cells.Item(426, 2).Value = 'Mártir De Cuilapan'
$ws.Cells.Item(439, 2).Value = 'Taxco De Alarcón'
$ws.Cells.Item(441, 2).Value = 'Técpan De Galeana'
$ws.Cells.Item(443, 2).Value = 'Tepecoacuilco De Trujano'
$ws.Cells.Item(445, 2).Value = 'Tixtla De Guerrero'
$ws.Cells.Item(449, 2).Value = 'Tlalixtaquilla De Maldonado'
$ws.Cells.Item(450, 2).Value = 'Tlapa De Comonfort'
$ws.Cells.Item(462, 2).Value = 'Agua Blanca De Iturbide'
$ws.Cells.Item(468, 2).Value = 'Atotonilco De Tula'
$ws.Cells.Item(469, 2).Value = 'Atotonilco El Grande'
$ws.Cells.Item(475, 2).Value = 'Cuautepec De Hinojosa'
$ws.Cells.Item(480, 2).Value = 'Huasca De Ocampo'
$ws.Cells.Item(484, 2).Value = 'Huejutla De Reyes'
$ws.Cells.Item(487, 2).Value = 'Jacala De Ledezma'
$ws.Cells.Item(493, 2).Value = 'Mineral De La Reforma'
$ws.Cells.Item(494, 2).Value = 'Mineral Del Chico'
$ws.Cells.Item(495, 2).Value = 'Mineral Del Monte'
$ws.Cells.Item(496, 2).Value = 'Mixquiahuala De Juárez'
$ws.Cells.Item(497, 2).Value = 'Molango De Escamilla'
$ws.Cells.Item(499, 2).Value = 'Nopala De Villagrán'
$ws.Cells.Item(500, 2).Value = 'Omitlán De Juárez'
$ws.Cells.Item(501, 2).Value = 'Pachuca De Soto'
$ws.Cells.Item(504, 2).Value = 'Progreso De Obregón'
$ws.Cells.Item(509, 2).Value = 'Santiago De Anaya'
$ws.Cells.Item(510, 2).Value = 'Santiago Tulantepec De Lugo Guerrero'
$ws.Cells.Item(514, 2).Value = 'Tenango De Doria'
$ws.Cells.Item(516, 2).Value = 'Tepehuacán De Guerrero'
$ws.Cells.Item(517, 2).Value = 'Tepeji Del Río De Ocampo'
$ws.Cells.Item(519, 2).Value = 'Tezontepec De Aldama'
$ws.Cells.Item(527, 2).Value = 'Tula De Allende'
$ws.Cells.Item(528, 2).Value = 'Tulancingo De Bravo'
$ws.Cells.Item(529, 2).Value = 'Villa De Tezontepec'
$ws.Cells.Item(533, 2).Value = 'Zacualtipán De Ángeles'
$ws.Cells.Item(538, 2).Value = 'Acatlán De Juárez'
$ws.Cells.Item(539, 2).Value = 'Ahualulco De Mercado'
$ws.Cells.Item(543, 2).Value = 'Atotonilco El Alto'
$ws.Cells.Item(545, 2).Value = 'Autlán De Navarro'
$ws.Cells.Item(560, 2).Value = 'Encarnación De Díaz'
$ws.Cells.Item(567, 2).Value = 'Huejuquilla El Alto'
$ws.Cells.Item(568, 2).Value = 'Ixtlahuacán De Los Membrillos'
$ws.Cells.Item(569, 2).Value = 'Ixtlahuacán Del Río'
$ws.Cells.Item(573, 2).Value = 'Jilotlán De Los Dolores'
$ws.Cells.Item(579, 2).Value = 'La Manzanilla De La Paz'
$ws.Cells.Item(580, 2).Value = 'Lagos De Moreno'
$ws.Cells.Item(587, 2).Value = 'Ojuelos De Jalisco'
$ws.Cells.Item(592, 2).Value = 'San Cristóbal De La Barranca'
$ws.Cells.Item(593, 2).Value = 'San Diego De Alejandría'
$ws.Cells.Item(595, 2).Value = 'San Juan De Los Lagos'
$ws.Cells.Item(597, 2).Value = 'San Martín De Bolaños'
$ws.Cells.Item(598, 2).Value = 'San Miguel El Alto'
$ws.Cells.Item(599, 2).Value = 'San Sebastián Del Oeste'
$ws.Cells.Item(600, 2).Value = 'Santa María De Los Ángeles'
$ws.Cells.Item(601, 2).Value = 'Santa María Del Oro'
$ws.Cells.Item(604, 2).Value = 'Talpa De Allende'
$ws.Cells.Item(605, 2).Value = 'Tamazula De Gordiano'
$ws.Cells.Item(608, 2).Value = 'Techaluta De Montenegro'
$ws.Cells.Item(611, 2).Value = 'Teocuitatlán De Corona'
$ws.Cells.Item(612, 2).Value = 'Tepatitlán De Morelos'
$ws.Cells.Item(615, 2).Value = 'Tizapán El Alto'
$ws.Cells.Item(616, 2).Value = 'Tlajomulco De Zúñiga'
$ws.Cells.Item(626, 2).Value = 'Unión De Tula'
$ws.Cells.Item(627, 2).Value = 'Valle De Guadalupe'
$ws.Cells.Item(628, 2).Value = 'Valle De Juárez'
$ws.Cells.Item(633, 2).Value = 'Yahualica De González Gallo'
$ws.Cells.Item(634, 2).Value = 'Zacoalco De Torres'
$ws.Cells.Item(637, 2).Value = 'Zapotitlán De Vadillo'
$ws.Cells.Item(638, 2).Value = 'Zapotlán Del Rey'
$ws.Cells.Item(639, 2).Value = 'Zapotlán El Grande'
$ws.Cells.Item(664, 2).Value = 'Coalcomán De Vázquez Pallares'
$ws.Cells.Item(666, 2).Value = 'Cojumatlán De Régules'
$ws.Cells.Item(732, 2).Value = 'Tiquicheo De Nicolás Romero'
$ws.Cells.Item(756, 2).Value = 'Coatlán Del Río'
$ws.Cells.Item(767, 2).Value = 'Puente De Ixtla'
$ws.Cells.Item(773, 2).Value = 'Tetela Del Volcán'
$ws.Cells.Item(775, 2).Value = 'Tlaltizapán De Zapata'
$ws.Cells.Item(783, 2).Value = 'Zacualpan De Amilpas'
$ws.Cells.Item(787, 2).Value = 'Amatlán De Cañas'
$ws.Cells.Item(788, 2).Value = 'Bahía De Banderas'
$ws.Cells.Item(790, 2).Value = 'Ixtlán Del Río'
$ws.Cells.Item(797, 2).Value = 'Santa María Del Oro'
$ws.Cells.Item(817, 2).Value = 'Mier Y Noriega'
$ws.Cells.Item(822, 2).Value = 'San Nicolás De Los Garza'
$ws.Cells.Item(829, 2).Value = 'Acatlán De Pérez Figueroa'
$ws.Cells.Item(836, 2).Value = 'Chalcatongo De Hidalgo'
$ws.Cells.Item(838, 2).Value = 'Coicoyán De Las Flores'
$ws.Cells.Item(839, 2).Value = 'Constancia Del Rosario'
$ws.Cells.Item(841, 2).Value = 'Guevea De Humboldt'
$ws.Cells.Item(842, 2).Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Cells.Item(843, 2).Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Cells.Item(844, 2).Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Cells.Item(845, 2).Value = 'Ixtlán De Juárez'
$ws.Cells.Item(846, 2).Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Cells.Item(854, 2).Value = 'Magdalena Yodocono De Porfirio Díaz'
$ws.Cells.Item(855, 2).Value = 'Mariscala De Juárez'
$ws.Cells.Item(856, 2).Value = 'Mártires De Tacubaya'
$ws.Cells.Item(858, 2).Value = 'Mazatlán Villa De Flores'
$ws.Cells.Item(860, 2).Value = 'Miahuatlán De Porfirio Díaz'
$ws.Cells.Item(864, 2).Value = 'Nejapa De Madero'
$ws.Cells.Item(865, 2).Value = 'Oaxaca De Juárez'
$ws.Cells.Item(866, 2).Value = 'Ocotlán De Morelos'
$ws.Cells.Item(867, 2).Value = 'Pinotepa De Don Luis'
$ws.Cells.Item(869, 2).Value = 'Putla Villa De Guerrero'
$ws.Cells.Item(870, 2).Value = 'Reforma De Pineda'
$ws.Cells.Item(881, 2).Value = 'San Antonino El Alto'
$ws.Cells.Item(895, 2).Value = 'San Francisco Del Mar'
$ws.Cells.Item(906, 2).Value = 'San José Del Progreso'
$ws.Cells.Item(914, 2).Value = 'San Juan Bautista Lo De Soto'
$ws.Cells.Item(951, 2).Value = 'San Mateo Del Mar'
$ws.Cells.Item(959, 2).Value = 'San Miguel Del Puerto'
$ws.Cells.Item(961, 2).Value = 'San Miguel El Grande'
$ws.Cells.Item(975, 2).Value = 'San Pablo Villa De Mitla'
$ws.Cells.Item(976, 2).Value = 'San Pedro El Alto'
$ws.Cells.Item(988, 2).Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Cells.Item(1006, 2).Value = 'Santa Cruz De Bravo'
$ws.Cells.Item(1010, 2).Value = 'Santa Cruz Tacache De Mina'
$ws.Cells.Item(1014, 2).Value = 'Santa Inés De Zaragoza'
$ws.Cells.Item(1020, 2).Value = 'Santa María Del Tule'
$ws.Cells.Item(1026, 2).Value = 'Santa María Jalapa Del Marqués'
$ws.Cells.Item(1072, 2).Value = 'Santo Domingo De Morelos'
$ws.Cells.Item(1086, 2).Value = 'Sitio De Xitlapehua'
$ws.Cells.Item(1087, 2).Value = 'Tamazulápam Del Espíritu Santo'
$ws.Cells.Item(1089, 2).Value = 'Tataltepec De Valdés'
$ws.Cells.Item(1090, 2).Value = 'Teotitlán De Flores Magón'
$ws.Cells.Item(1092, 2).Value = 'Tepelmeme Villa De Morelos'
$ws.Cells.Item(1094, 2).Value = 'Tlacolula De Matamoros'
$ws.Cells.Item(1095, 2).Value = 'Totontepec Villa De Morelos'
$ws.Cells.Item(1098, 2).Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Cells.Item(1099, 2).Value = 'Villa De Zaachila'
$ws.Cells.Item(1101, 2).Value = 'Villa Sola De Vega'
$ws.Cells.Item(1102, 2).Value = 'Villa Talea De Castro'
$ws.Cells.Item(1106, 2).Value = 'Zimatlán De Álvarez'
$ws.Cells.Item(1127, 2).Value = 'Chalchicomula De Sesma'
$ws.Cells.Item(1136, 2).Value = 'Chila De La Sal'
$ws.Cells.Item(1144, 2).Value = 'Cuayuca De Andrade'
$ws.Cells.Item(1145, 2).Value = 'Cuetzalan Del Progreso'
$ws.Cells.Item(1158, 2).Value = 'Huehuetlán El Chico'
$ws.Cells.Item(1159, 2).Value = 'Huehuetlán El Grande'
$ws.Cells.Item(1164, 2).Value = 'Ixcamilpa De Guerrero'
$ws.Cells.Item(1167, 2).Value = 'Izúcar De Matamoros'
$ws.Cells.Item(1178, 2).Value = 'Los Reyes De Juárez'
$ws.Cells.Item(1188, 2).Value = 'Palmar De Bravo'
$ws.Cells.Item(1208, 2).Value = 'San Nicolás De Los Ranchos'
$ws.Cells.Item(1212, 2).Value = 'San Salvador El Verde'
$ws.Cells.Item(1218, 2).Value = 'Tecali De Herrera'
$ws.Cells.Item(1224, 2).Value = 'Tepanco De López'
$ws.Cells.Item(1225, 2).Value = 'Tepango De Rodríguez'
$ws.Cells.Item(1226, 2).Value = 'Tepatlaxco De Hidalgo'
$ws.Cells.Item(1230, 2).Value = 'Tepexi De Rodríguez'
$ws.Cells.Item(1232, 2).Value = 'Tetela De Ocampo'
$ws.Cells.Item(1233, 2).Value = 'Teteles De Avila Castillo'
$ws.Cells.Item(1238, 2).Value = 'Tlacotepec De Benito Juárez'
$ws.Cells.Item(1251, 2).Value = 'Xayacatlán De Bravo'
$ws.Cells.Item(1256, 2).Value = 'Xochitlán De Vicente Suárez'
$ws.Cells.Item(1267, 2).Value = 'Amealco De Bonfil'
$ws.Cells.Item(1269, 2).Value = 'Cadereyta De Montes'
$ws.Cells.Item(1274, 2).Value = 'Jalpan De Serra'
$ws.Cells.Item(1275, 2).Value = 'Landa De Matamoros'
$ws.Cells.Item(1278, 2).Value = 'Pinal De Amoles'
$ws.Cells.Item(1281, 2).Value = 'San Juan Del Río'
$ws.Cells.Item(1293, 2).Value = 'Armadillo De Los Infante'
$ws.Cells.Item(1294, 2).Value = 'Axtla De Terrazas'
$ws.Cells.Item(1300, 2).Value = 'Ciudad Del Maíz'
$ws.Cells.Item(1309, 2).Value = 'Mexquitic De Carmona'
$ws.Cells.Item(1315, 2).Value = 'San Ciro De Acosta'
$ws.Cells.Item(1321, 2).Value = 'Santa María Del Río'
$ws.Cells.Item(1323, 2).Value = 'Soledad De Graciano Sánchez'
$ws.Cells.Item(1330, 2).Value = 'Tanquián De Escobedo'
$ws.Cells.Item(1334, 2).Value = 'Villa De Arista'
$ws.Cells.Item(1335, 2).Value = 'Villa De Arriaga'
$ws.Cells.Item(1336, 2).Value = 'Villa De Guadalupe'
$ws.Cells.Item(1337, 2).Value = 'Villa De La Paz'
$ws.Cells.Item(1338, 2).Value = 'Villa De Ramos'
$ws.Cells.Item(1339, 2).Value = 'Villa De Reyes'
$ws.Cells.Item(1369, 2).Value = 'Nacozari De García'
$ws.Cells.Item(1384, 2).Value = 'Jalpa De Méndez'
$ws.Cells.Item(1422, 2).Value = 'Soto La Marina'
$ws.Cells.Item(1436, 2).Value = 'Mazatecochco De José María Morelos'
$ws.Cells.Item(1437, 2).Value = 'Nanacamilpa De Mariano Arista'
$ws.Cells.Item(1440, 2).Value = 'Papalotla De Xicohténcatl'
$ws.Cells.Item(1441, 2).Value = 'San Pablo Del Monte'
$ws.Cells.Item(1445, 2).Value = 'Tepetitla De Lardizábal'
$ws.Cells.Item(1448, 2).Value = 'Tetla De La Solidaridad'
$ws.Cells.Item(1465, 2).Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Cells.Item(1469, 2).Value = 'Amatlán De Los Reyes'
$ws.Cells.Item(1479, 2).Value = 'Boca Del Río'
$ws.Cells.Item(1481, 2).Value = 'Camarón De Tejeda'
$ws.Cells.Item(1485, 2).Value = 'Castillo De Teayo'
$ws.Cells.Item(1487, 2).Value = 'Cazones De Herrera'
$ws.Cells.Item(1504, 2).Value = 'Cosamaloapan De Carpio'
$ws.Cells.Item(1505, 2).Value = 'Cosautlán De Carvajal'
$ws.Cells.Item(1521, 2).Value = 'Hueyapan De Ocampo'
$ws.Cells.Item(1522, 2).Value = 'Ignacio De La Llave'
$ws.Cells.Item(1526, 2).Value = 'Ixhuacán De Los Reyes'
$ws.Cells.Item(1527, 2).Value = 'Ixhuatlán De Madero'
$ws.Cells.Item(1528, 2).Value = 'Ixhuatlán Del Sureste'
$ws.Cells.Item(1538, 2).Value = 'Juchique De Ferrer'
$ws.Cells.Item(1542, 2).Value = 'Lerdo De Tejada'
$ws.Cells.Item(1545, 2).Value = 'Martínez De La Torre'
$ws.Cells.Item(1548, 2).Value = 'Medellín De Bravo'
$ws.Cells.Item(1552, 2).Value = 'Mixtla De Altamirano'
$ws.Cells.Item(1554, 2).Value = 'Nanchital De Lázaro Cárdenas Del Río'
$ws.Cells.Item(1561, 2).Value = 'Ozuluama De Mascareñas'
$ws.Cells.Item(1565, 2).Value = 'Paso De Ovejas'
$ws.Cells.Item(1566, 2).Value = 'Paso Del Macho'
$ws.Cells.Item(1570, 2).Value = 'Poza Rica De Hidalgo'
$ws.Cells.Item(1578, 2).Value = 'Sayula De Alemán'
$ws.Cells.Item(1581, 2).Value = 'Soledad De Doblado'
$ws.Cells.Item(1617, 2).Value = 'Vega De Alatorre'
$ws.Cells.Item(1628, 2).Value = 'Zontecomatlán De López Y Fuentes'
$ws.Cells.Item(1639, 2).Value = 'Cañitas De Felipe Pescador'
$ws.Cells.Item(1641, 2).Value = 'Concepción Del Oro'
$ws.Cells.Item(1652, 2).Value = 'Jiménez Del Teul'
$ws.Cells.Item(1658, 2).Value = 'Mezquital Del Oro'
$ws.Cells.Item(1662, 2).Value = 'Moyahua De Estrada'
$ws.Cells.Item(1663, 2).Value = 'Nochistlán De Mejía'
$ws.Cells.Item(1664, 2).Value = 'Noria De Ángeles'
$ws.Cells.Item(1674, 2).Value = 'Teúl De González Ortega'
$ws.Cells.Item(1675, 2).Value = 'Tlaltenango De Sánchez Román'
$ws.Cells.Item(1678, 2).Value = 'Villa De Cos'

# --- Fix floating point total recalculation for Guanajuato section (row 382) ---
$ws.Cells.Item(382, 4).Value = 0.09636532006019216

# --- Remove trailing footnote/metadata rows (1687-1691) ---
$ws.Range("A1687:A1691").EntireRow.Delete()

